# CONFIG_TWEAKS.docx tweaks before jumping to 2.0.7.2
#
# 1) Add three more Configuration.h list items right after
#    "#define Z_SAFE_HOMING", matching the existing ListParagraph /
#    numId=2 bullet formatting used by the rest of that list.
# 2) Unhide the built-in "Default Paragraph Font" character style
#    (drop its semiHidden flag).

$d = $word.ActiveDocument

# --- 1) Insert the three new #define bullets ---------------------------
$anchorText = "#define Z_SAFE_HOMING"
$newItems = @(
    "#define Z_MIN_POS -5",
    "//#define BLTOUCH_HS_MODE",
    "#define E0_CURRENT      1000"
)

$r = $d.Content
$found = $r.Find.Execute($anchorText, $false, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)

if ($found) {
    # Collapse to the end of the matched text (just before its paragraph
    # mark) so each InsertParagraphAfter() splits off a new list item that
    # inherits the "ListParagraph" style / numPr of the anchor paragraph.
    $r.Collapse(0)
    foreach ($item in $newItems) {
        $r.InsertParagraphAfter()
        $r.Collapse(0)
        $r.MoveStart(1, 1)
        $r.InsertBefore($item)
        $r.Collapse(0)
    }
}

# --- 2) Unhide the DefaultParagraphFont style ---------------------------
try {
    $style = $d.Styles("Default Paragraph Font")
    $style.Hidden = $false
} catch {
    # Some hosts treat this built-in style's visibility as read-only;
    # nothing more we can do through the object model in that case.
}
